# Ghostly Spymaster Design Doc - Iteration 1 edit
#
# The only real wording change in this revision lives in the "Environment"
# bullet about tile types. The paragraph that used to read as one big run of
# text:
#
#   "Because it is a top down 2D game we will be using Art tiles to define
#    what is ground, walls or environmental props. Wall tiles will have
#    hitboxes to break cone of vision and glass tiles will also have
#    hitboxes but does not break cone of vision, floor tiles are passable
#    by players and NPCs, Invisible hitbox tiles to restrict the player in
#    going to certain places."
#
# gets reformatted into a short line break after "...environmental props."
# plus three separate "Wall tiles -", "Floor tiles -", "Invisible hitbox
# tiles -" bullet-style lines (each its own paragraph, reusing the same
# ListParagraph / 360-twip indent formatting).

$d = $word.ActiveDocument

# 1) Break "...environmental props." onto its own line and start labelling
#    the wall-tile bullet with a dash: "Wall tiles - will have hitboxes..."
$d.Content.Find.Execute(
    "environmental props. Wall tiles", $true, $false, $false, $false, $false,
    $true, 1, $false, "environmental props. ^lWall tiles -", 2) | Out-Null

# 2) Word re-anchors its "last edit" (_GoBack) bookmark at the point of the
#    most recent change, i.e. right after "...environmental props. " and
#    before the new line break.
$anchor = $d.Content
$anchor.Find.Execute(
    "environmental props. ", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($anchor.End, $anchor.End))

# 3) Split the floor-tile sentence into its own bullet line: "Floor tiles -
#    are passable by players and NPCs."
$d.Content.Find.Execute(
    "cone of vision, floor tiles", $true, $false, $false, $false, $false,
    $true, 1, $false, "cone of vision. ^pFloor tiles -", 2) | Out-Null

# 4) Split the invisible-hitbox-tile sentence into its own bullet line:
#    "Invisible hitbox tiles - to restrict the player in going to certain
#    places."
$d.Content.Find.Execute(
    "NPCs, Invisible hitbox tiles to restrict", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "NPCs. ^pInvisible hitbox tiles - to restrict", 2) | Out-Null

Write-Output "Environment bullet split into Wall/Floor/Invisible-hitbox tile lines."
